# Updates to closed-loop recycling processes and excluding landfill and
# open-loop recycling from analysis.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new parameter row for closed-loop reprocessing of flat glass,
# directly below the existing "Recycling & Credits" rows (old row 41
# "Stillage & Logistics" and everything after it shifts down by one).
$ws.Rows.Item(41).Insert()

$ws.Range("A41").Value = "Recycling & Credits"
$ws.Range("B41").Value = "FLOAT_GLASS_REPROCESSING"
$ws.Range("C41").Value = 0.52
$ws.Range("D41").Value = "Embodied Carbon: Reprocessing Flat Glass (kgCO2e/kg)"

# Revised repurpose emission factors (closed-loop reuse updates).
$ws.Range("C28").Value = 1      # REPURPOSE_HEAVY_KGCO2_PER_M2   : 3 -> 1
$ws.Range("C29").Value = 0.5    # REPURPOSE_LIGHT_KGCO2_PER_M2   : 1 -> 0.5
$ws.Range("C30").Value = 0.75   # REPURPOSE_MEDIUM_KGCO2_PER_M2  : 2 -> 0.75
